$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 / Row 18 get a new E column value (no special style)
$ws.Range("E17").Value = "34653t"
$ws.Range("E18").Value = "re"

# New rows 19-28 with values in column E
$ws.Range("E19").Value = "yer"
$ws.Range("E20").Value = "ter"
$ws.Range("E21").Value = "ter"
$ws.Range("E22").Value = "ter"
$ws.Range("E23").Value = "ter"
$ws.Range("E24").Value = "ter"
$ws.Range("E25").Value = "`rter"
$ws.Range("E26").Value = "ter"
$ws.Range("E27").Value = "erywer5"
$ws.Range("E28").Value = "yer"

# Update selection to match the diff
$ws.Range("E28").Select()
